$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet/tab (was "Acc_Repayment1").
$ws.Activate()

# Insert a new (blank) column before column N, shifting
# Late / heading / Disbursement one column to the right (N->O, O->P, P->Q).
$ws.Columns("N:N").Insert()

# The newly inserted column takes on the same display width as column M.
$ws.Columns("N:N").ColumnWidth = 9.83

# Restore the cursor/selection to the cell it ended up on after the edit.
$null = $ws.Range("R11").Select()
